# Insert a new weekly price record as row 36 (pushing the existing
# historical rows 36-60 down to 37-61), matching the author's commit
# "Fruta / hortaliza, semanal" (weekly fruit/vegetable price update).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 36; rows 36..60 shift down to 37..61.
$ws.Rows.Item(36).Insert()

# Populate the new row 36 with this week's Ají (Americana) price record.
$ws.Range("A36").Value = 7
$ws.Range("B36").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C36").Value = "Ñuble"
$ws.Range("D36").Value = 44586
$ws.Range("E36").Value = 16
$ws.Range("F36").Value = 100112021
$ws.Range("G36").Value = "Ají"
$ws.Range("H36").Value = "Americana (o)"
$ws.Range("I36").Value = "Primera"
$ws.Range("J36").Value = 60
$ws.Range("K36").Value = 13000
$ws.Range("L36").Value = 14000
$ws.Range("M36").Value = 13500
$ws.Range("N36").Value = "`$/caja 15 kilos"
$ws.Range("O36").Value = "Región del Maule"
$ws.Range("P36").Value = 900
$ws.Range("Q36").Value = 15
$ws.Range("R36").Value = "Hortaliza"
